# Florida extended workbook - Wisconsin-data commit's incidental edits on Sheet1:
#  1. Header I1 renamed: "Voter Turnout Increase (VTI)" -> "Voter Turnout Change (VTI)"
#  2. Header M1 cleared (the "Gained Significant Republican Votes?" column header removed)
#  3. Column K (rows 2-68) formula direction flipped: $B-$E  ->  $E-$B
#  4. Selection moved from M3 to the I1:L1048576 column-band (I1 active)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Voter Turnout Increase (VTI)" header to "Voter Turnout Change (VTI)"
$ws.Range("I1").Value = "Voter Turnout Change (VTI)"

# 2. Clear the old "Gained Significant Republican Votes?" header cell
$ws.Range("M1").ClearContents()

# 3. Flip the K column's formula: now Voter Turnout Change = 2016 value - 2012 value (was 2012 - 2016)
#    Assigning one formula to the whole range lets Excel auto-adjust the relative references
#    per row, and it rebuilds the existing shared-formula groups.
$ws.Range("K2:K68").Formula = "=`$E2-`$B2"

# 4. Update the sheet's selection/view to the I:L column band with I1 as the active cell
[void]$ws.Range("I1:L1048576").Select()
